$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 content: replace the old joke answers with the new ones,
# and fill in the previously-empty banniere (C4) cell with an image URL.
$ws.Range("A4").Value = "Wthl we are on it !"
$ws.Range("B4").Value = "Oh yeah homie we are on it !"
$ws.Range("C4").Value = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcSGRJtbkKjWAfMGYQQ652F1xxK-JRDoiZ1Znw&s"
